$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the serial numbers in column A (rows 11-37): after deleting some
# earlier rows the numbering kept the old, now-skipped values. Renumber
# them so the sequence is contiguous again (10, 11, 12, ... 36).
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23
$ws.Range("A25").Value = 24
$ws.Range("A26").Value = 25
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29
$ws.Range("A31").Value = 30
$ws.Range("A32").Value = 31
$ws.Range("A33").Value = 32
$ws.Range("A34").Value = 33
$ws.Range("A35").Value = 34
$ws.Range("A36").Value = 35
$ws.Range("A37").Value = 36

# A37 was left with the mismatched (bold-ish Calibri) style from before the
# rows above it were cleaned up; match it to the rest of the column (A2:A36)
# by copying that formatting over from the cell right above it.
$ws.Range("A36").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-select the numbering column to review the fix, scrolled down to the
# last rows that were just edited.
$ws.Range("A2:A37").Select()
